$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TestCases (B2) and Instance (D2) values - BAL case split
$ws.Range("B2").Value = "42-45"
$ws.Range("D2").Value = "Automation1"

# Update the active selection to D3
$ws.Range("D3").Select()
